$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 439, shifting existing rows 439:515 down to 440:516.
$ws.Rows(439).Insert()

# Populate the newly inserted row 439 with the new record's values.
$ws.Range("A439").Value = 6
$ws.Range("B439").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C439").Value = "Metropolitana"
$ws.Range("D439").Value = 44798
$ws.Range("E439").Value = 13
$ws.Range("F439").Value = 100112052
$ws.Range("G439").Value = "Albahaca"
$ws.Range("H439").Value = "Sin especificar"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 300
$ws.Range("K439").Value = 3500
$ws.Range("L439").Value = 4000
$ws.Range("M439").Value = 3717
$ws.Range("N439").Value = "`$/paquete"
$ws.Range("O439").Value = "Limache"
$ws.Range("P439").Value = 3717
$ws.Range("Q439").Value = 1
$ws.Range("R439").Value = "Hortaliza"
